# 自动更新Excel文件 - 2025-11-02 23:11:34
# For each data row, decrement the "剩余" (remaining) count in column E by 1.
# When remaining reaches 1 (i.e. would hit 0), reset it back to the full
# cycle length stored in column D ("总天") and roll the "开始时间" (start
# date) in column F forward to the new cycle start date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStartDate = 20251103

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 2
}

for ($row = 2; $row -le $lastRow; $row++) {
    $totalDaysCell = $ws.Cells.Item($row, 4)   # column D - 总天
    $remainingCell = $ws.Cells.Item($row, 5)   # column E - 剩余
    $startDateCell = $ws.Cells.Item($row, 6)   # column F - 开始时间

    $totalDays = $totalDaysCell.Value2
    $remaining = $remainingCell.Value2
    $startDate = $startDateCell.Value2

    if ($remaining -eq $null -or $totalDays -eq $null) {
        continue
    }

    # The start date must be a well-formed 8-digit YYYYMMDD value; rows with
    # a malformed date (e.g. a stray extra digit) are left untouched, same
    # as the source process that produced this update skipped them.
    $startDateText = [string]$startDate
    if ($startDateText.Length -ne 8) {
        continue
    }

    if ($remaining -eq 1) {
        $remainingCell.Value2 = $totalDays
        $startDateCell.Value2 = $newStartDate
    }
    else {
        $remainingCell.Value2 = $remaining - 1
    }
}
